$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 table: change the table style (tableStyleId) on the graphic-frame
#    table, e.g. via Table.ApplyStyle(styleId, isTableStyleId).
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{30868731-10E6-4122-A7D8-F27734FAEC68}", $true)
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the deck's live theme (the one bound to the slide master /
#    design) from the "Integral" / Red Violet palette over to the
#    "Office Theme" palette, matching the swap recorded in the diff.
# ---------------------------------------------------------------------------
# Index order exposed by ThemeColorScheme / Master.ColorScheme:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB is the standard OLE_COLOR (BGR-packed) integer used by VBA's RGB().
$officeColorsBgr = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $officeColorsBgr[$i - 1]
}
